$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date placeholder text on the
#    slide master and every slide layout (15.09.2025 -> 19.09.2025).
# ---------------------------------------------------------------------------
$oldDate = "15.09.2025"
$newDate = "19.09.2025"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Datumsplatzhalter*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout (custom layout) hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) Rewrite the intro text body on slide 6 ("Header" / "abschnitte") with
#    the new "Bplaced - kostenloser Hostingservice" copy.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Build the first paragraph as three separate runs (so later per-run
# formatting stays isolated) and leave a second, empty paragraph behind.
$tr.Text = "Bplaced"
$run2 = $tr.InsertAfter(" " + [char]0x2013 + " kostenloser ")
$run3 = $run2.InsertAfter("Hostingservice" + [char]13)
